$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update summary figures at top of the statement
#    VALOR MORA total increases by one extra period (40000)
#    Cant. Periodos increases from 40 to 41
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1618666
$ws.Range("F13").Value = 41

# ---------------------------------------------------------------------------
# 2) Rebuild the period detail rows (16-56).
#    Previously rows 16-55 listed periods 2507 down to 2204 (descending),
#    with the final row (2204) carrying a partial value of 18666 and all
#    other periods carrying 40000.
#    Now rows 16-56 list periods 2204 up to 2508 (ascending) - one new
#    period (2508) was appended - and the partial value (18666) stays on
#    period 2204, which is now the first row instead of the last.
# ---------------------------------------------------------------------------

# 2a) Capture the special "last row" formatting (currently row 55) before
#     we touch anything, and stamp it onto the new row 56.
$ws.Range("B55:J55").Copy()
$ws.Range("B56:J56").PasteSpecial(-4122)

# 2b) Convert old row 55 to the regular detail-row formatting (copy from
#     row 54, which already uses the standard style).
$ws.Range("B54:J54").Copy()
$ws.Range("B55:J55").PasteSpecial(-4122)

# 2c) Write out B/C/D/G (constant across every detail row) for the new row.
$ws.Range("B56").Value = $ws.Range("B16").Text
$ws.Range("C56").Value = $ws.Range("C16").Text
$ws.Range("D56").Value = $ws.Range("D16").Text
$ws.Range("G56").Value = 1000000

# 2d) Build the ascending period list 2204 .. 2508 and assign E (period)
#     and F (value) for rows 16 through 56.
$periods = @()
foreach ($m in 4..12) { $periods += ("22{0:D2}" -f $m) }
foreach ($m in 1..12) { $periods += ("23{0:D2}" -f $m) }
foreach ($m in 1..12) { $periods += ("24{0:D2}" -f $m) }
foreach ($m in 1..8)  { $periods += ("25{0:D2}" -f $m) }

$row = 16
foreach ($period in $periods) {
    $ws.Range("E$row").Value = $period
    if ($period -eq "2204") {
        $ws.Range("F$row").Value = 18666
    } else {
        $ws.Range("F$row").Value = 40000
    }
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3) Shift the signature block down by one row: the blank underline row
#    (old row 60) becomes row 61, and the NOMBRE/FIRMA label row (old row
#    61) becomes row 62. This makes room for the extra detail row above.
# ---------------------------------------------------------------------------
$row61_B = $ws.Range("B61").Text
$row61_H = $ws.Range("H61").Text
$row60_B = $ws.Range("B60").Text
$row60_H = $ws.Range("H60").Text

# 3a) Move the NOMBRE/FIRMA label row to row 62.
$ws.Range("B61:C61").Copy()
$ws.Range("B62:C62").PasteSpecial(-4122)
$ws.Range("H61:J61").Copy()
$ws.Range("H62:J62").PasteSpecial(-4122)
$ws.Range("B62:C62").Merge()
$ws.Range("H62:J62").Merge()
$ws.Range("B62").Value = $row61_B
$ws.Range("H62").Value = $row61_H

# 3b) Move the blank underline row to row 61.
$ws.Range("B60:C60").Copy()
$ws.Range("B61:C61").PasteSpecial(-4122)
$ws.Range("H60:J60").Copy()
$ws.Range("H61:J61").PasteSpecial(-4122)
$ws.Range("B61:C61").Merge()
$ws.Range("H61:J61").Merge()
$ws.Range("B61").Value = $row60_B
$ws.Range("H61").Value = $row60_H

# 3c) Clear out the now-unused row 60 (content, formatting and merges).
$ws.Range("B60:C60").UnMerge()
$ws.Range("H60:J60").UnMerge()
$ws.Range("B60:J60").Clear()
